$d = $word.ActiveDocument

# 1. Simple in-place text swap: codigo_comercio -> central_bank_trade_code
$null = $d.Content.Find.Execute("codigo_comercio", $false, $false, $false, $false, $false, $true, 1, $false, "central_bank_trade_code", 2)

# 2. Locate the paragraph index whose range contains the first match of $searchText.
function Get-ParaIndexContainingText {
    param([string]$searchText)
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        return -1
    }
    $target = $rng.Start
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($target -ge $p.Range.Start -and $target -lt $p.Range.End) {
            return $i
        }
    }
    return -1
}

$startIdx = Get-ParaIndexContainingText "Por favor {action_currency_1}"
$endIdx   = Get-ParaIndexContainingText "account_bank_currency_2"

$pStart = $d.Paragraphs.Item($startIdx)
$pEnd   = $d.Paragraphs.Item($endIdx)

# Range spanning both "Por favor ..." paragraphs, including the trailing paragraph mark.
$targetRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$newParasXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="120" w:line="240" w:lineRule="auto" /><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>CUENTA PARA RECIBIR {</w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>abonar_currency</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>}:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="120" w:line="240" w:lineRule="auto" /><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr></w:pPr><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>Número</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>: {</w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>abonar_account_number</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="240" w:line="240" w:lineRule="auto" /><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>Banco: {</w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>abonar_bank_name</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="120" w:line="240" w:lineRule="auto" /><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>CUENTA PARA PAGAR {</w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>cargar_currency</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>}:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="120" w:line="240" w:lineRule="auto" /><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>Número: {</w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>cargar_account_number</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="240" w:line="240" w:lineRule="auto" /><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>Banco: {</w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>cargar_bank_name</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Manrope" w:hAnsi="Manrope" /><w:sz w:val="20" /><w:szCs w:val="20" /><w:lang w:val="es-CL" /></w:rPr><w:t>}</w:t></w:r></w:p>
'@

$targetRange.InsertXML($newParasXml)
